$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 new rows directly below the header (row 1), pushing the existing
#    20 data rows (originally rows 2-21, timestamps 0-1900) down to rows 5-24.
$ws.Rows("2:4").Insert()

# Excel's row-insert copies formatting down from the row above (the bold header),
# so strip that back off to keep the new rows looking like plain data rows.
$ws.Range("A2:H4").ClearFormats()

# 2. Fill the 3 newly inserted rows with their sensor readings (label + ax..gz).
#    Timestamps are reassigned below in step 4 for the whole A column.
$ws.Range("B2").Value = "falling"
$ws.Range("C2").Value = -0.2407464981079101
$ws.Range("D2").Value = 0.5433270186185837
$ws.Range("E2").Value = -0.3658644706010817
$ws.Range("F2").Value = 0.0885754600167274
$ws.Range("G2").Value = -0.0678060427308082
$ws.Range("H2").Value = -0.0534507073462009

$ws.Range("B3").Value = "falling"
$ws.Range("C3").Value = -0.2030138969421392
$ws.Range("D3").Value = 0.5583634674549103
$ws.Range("E3").Value = -0.1664191037416463
$ws.Range("F3").Value = 0.0888808965682983
$ws.Range("G3").Value = 0.1319468915462494
$ws.Range("H3").Value = -0.0287106670439243

$ws.Range("B4").Value = "falling"
$ws.Range("C4").Value = -0.3572314977645862
$ws.Range("D4").Value = 0.5254133790731429
$ws.Range("E4").Value = -0.5817861706018449
$ws.Range("F4").Value = 0.3381139039993286
$ws.Range("G4").Value = 0.00534507073462
$ws.Range("H4").Value = 0.1145372316241264

# 3. Append 7 more rows (rows 25-31) at the bottom with new sensor readings.
$ws.Range("B25").Value = "falling"
$ws.Range("C25").Value = -0.0122048854827883
$ws.Range("D25").Value = 0.3785421848297117
$ws.Range("E25").Value = -0.2985985279083258
$ws.Range("F25").Value = 0.1823432743549347
$ws.Range("G25").Value = 0.1472185254096985
$ws.Range("H25").Value = -0.3480404615402221

$ws.Range("B26").Value = "falling"
$ws.Range("C26").Value = 0.5578445792198222
$ws.Range("D26").Value = 0.3809743523597721
$ws.Range("E26").Value = 0.0945302546024358
$ws.Range("F26").Value = 0.5630650520324707
$ws.Range("G26").Value = 0.6291912198066711
$ws.Range("H26").Value = -0.8999572396278381

$ws.Range("B27").Value = "falling"
$ws.Range("C27").Value = 1.200664520263672
$ws.Range("D27").Value = 0.5283758044242859
$ws.Range("E27").Value = 0.6270142197608946
$ws.Range("F27").Value = -0.0210748501121997
$ws.Range("G27").Value = -0.0826195254921913
$ws.Range("H27").Value = -0.0035124751739203

$ws.Range("B28").Value = "falling"
$ws.Range("C28").Value = 0.02087068557739279
$ws.Range("D28").Value = 0.6865898966789239
$ws.Range("E28").Value = -0.5058017373085018
$ws.Range("F28").Value = -0.0372627787292003
$ws.Range("G28").Value = 0.0313068442046642
$ws.Range("H28").Value = -0.0632245540618896

$ws.Range("B29").Value = "falling"
$ws.Range("C29").Value = 0.09223079681396533
$ws.Range("D29").Value = 0.511084794998169
$ws.Range("E29").Value = -0.3808159828186036
$ws.Range("F29").Value = 0.0158824957907199
$ws.Range("G29").Value = 0.0445931628346443
$ws.Range("H29").Value = 0.0317649915814399

$ws.Range("B30").Value = "falling"
$ws.Range("C30").Value = 0.2094589471817022
$ws.Range("D30").Value = 0.502252608537674
$ws.Range("E30").Value = -0.392595499753952
$ws.Range("F30").Value = 0.0148134818300604
$ws.Range("G30").Value = -0.0259617734700441
$ws.Range("H30").Value = 0.0120645882561802

$ws.Range("B31").Value = "falling"
$ws.Range("C31").Value = 0.01102042198181016
$ws.Range("D31").Value = 0.4913336634635925
$ws.Range("E31").Value = -0.3189654350280756
$ws.Range("F31").Value = 0.0215329993516206
$ws.Range("G31").Value = 0.0319177098572254
$ws.Range("H31").Value = 0.0093156946823

# 4. The "timestamp" column (A) is a simple regenerated index: 0, 100, 200, ...
#    for every one of the 30 data rows now in the sheet (rows 2-31), independent
#    of which sensor-reading row ended up at that position.
for ($i = 0; $i -lt 30; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $i * 100
}
